# Apply the row-content reshuffle + trial_total increment described by the diff.
# Columns affected per data row (rows 2..41): F, H, I, K, L, M, N, O, P, Q, R, S, T, U, V
# F (trial_total) always increases by 81 for every row, in place.
# All the other listed columns get the content that a *different* source row
# (same column) held before the edit - i.e. the 40 data rows are permuted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> source row mapping (destination row = key, value = row whose old content is copied in)
$map = @{
    2  = 17
    3  = 37
    4  = 11
    5  = 34
    6  = 19
    7  = 40
    8  = 9
    9  = 23
    10 = 8
    11 = 15
    12 = 21
    13 = 41
    14 = 32
    15 = 10
    16 = 4
    17 = 6
    18 = 22
    19 = 30
    20 = 13
    21 = 3
    22 = 18
    23 = 38
    24 = 25
    25 = 28
    26 = 26
    27 = 14
    28 = 16
    29 = 5
    30 = 20
    31 = 36
    32 = 2
    33 = 35
    34 = 12
    35 = 24
    36 = 33
    37 = 29
    38 = 39
    39 = 27
    40 = 31
    41 = 7
}

# Columns (by number) that get copied from the source row, besides F.
# 8=H, 9=I, 11=K, 12=L, 13=M, 14=N, 15=O, 16=P, 17=Q, 18=R, 19=S, 20=T, 21=U, 22=V
$cols = @(8, 9, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22)

# First, snapshot the "before" state of every relevant cell for all 40 data rows,
# since several rows act as the source for other rows (we must read before overwriting).
$snapshot = @{}
for ($r = 2; $r -le 41; $r++) {
    $rowData = @{}
    $rowData[6] = $ws.Cells.Item($r, 6).Value()
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Cells.Item($r, $c).Value()
    }
    $snapshot[$r] = $rowData
}

# Now write the new values: F = old F + 81 (in place); other columns copied from source row snapshot.
for ($r = 2; $r -le 41; $r++) {
    $src = $map[$r]
    $srcData = $snapshot[$src]
    $ownData = $snapshot[$r]

    $ws.Cells.Item($r, 6).Value = $ownData[6] + 81

    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value = $srcData[$c]
    }
}
